$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width tweaks --------------------------------------------------
# (Target widths, as stored in the xlsx <col> elements, are snapped by the
# host's column-width<->pixel model; the input values below are chosen so the
# round-tripped width lands on the nearest achievable value to the target.)
$ws.Columns.Item(1).ColumnWidth = 19.142857142857142   # -> ~19.86 (target 19.9140625)
$ws.Columns.Item(7).ColumnWidth = 20.857142857142858   # -> ~21.57 (target 21.58203125)
$ws.Columns.Item(9).ColumnWidth = 21.142857142857142   # -> ~21.86 (target 21.83203125)

# --- New row 8: Interface / Manager ----------------------------------------
$ws.Range("A8").Value = "Interface"
$ws.Range("B8").Value = "cpp"
$ws.Range("C8").Value = "frame"
$ws.Range("F8").Value = "dll"
$ws.Range("G8").Value = "Manager"

# --- New row 7: worldBuilder-Minecraft (rebuilt) ---------------------------
$ws.Range("A7").Value = "worldBuilder-Minecraft"
$ws.Range("B7").Value = "py"
$ws.Range("E7").Value = "Minecraft"
$ws.Range("F7").Value = "pyd"
$ws.Range("G7").Value = "worldBuilder-Minecraft"
$ws.Range("H7").Value = "a1b2c1"

# --- New row 9: githubExplorer ----------------------------------------------
$ws.Range("A9").Value = "githubExplorer"
$ws.Range("B9").Value = "py"
$ws.Range("E9").Value = "API"
$ws.Range("F9").Value = "pyd"
$ws.Range("G9").Value = "githubExplorer"

# --- Back-fill the TYPE column for the new rows -----------------------------
$ws.Range("C7").Value = "addon"
$ws.Range("C9").Value = "addon module"

# --- New header cell (LOGIC PATH note) --------------------------------------
$ws.Range("I1").Value = "module只返回数据不干活"

# --- Selection / view -------------------------------------------------------
$ws.Range("I14").Select() | Out-Null
